$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold values that look numeric/percentage but must stay as
# literal text (matching their original "inlineStr" cell type). Setting
# NumberFormat to "@" (Text) on each cell individually before writing the
# value prevents the Excel COM layer from auto-converting the string into a
# number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.56%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.65%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.212"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.53%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07683"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.49%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.291"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.54%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.719"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.42%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9315"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.62%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.37%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1283"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.09%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.53%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09121"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.45%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04218"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.12%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.03%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001279"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.00%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005889"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.42%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.13%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.637"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "13.12%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1353"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.87%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2720"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.52%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04023"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.45%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001267"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.07%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004104"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.40%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.07%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02547"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05308"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.95%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007833"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.17%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1311"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.61%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006652"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.82%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002054"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008105"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.17%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.06%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006787"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.22%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2254"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "230.78%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.09%"
